# Update vm_pu.xlsx res_bus values for the 380 kV case (rows 2-25, cols B-F & I-N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032788237026491
$ws.Range("D2").Value = 1.033753609314567
$ws.Range("E2").Value = 1.040933865598175
$ws.Range("F2").Value = 1.04891546687882
$ws.Range("I2").Value = 1.031739051617225
$ws.Range("J2").Value = 1.037916532413789
$ws.Range("K2").Value = 1.036555050890173
$ws.Range("L2").Value = 1.043714808676399
$ws.Range("M2").Value = 1.051673981004004
$ws.Range("N2").Value = 1.039390492940507
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033772418484684
$ws.Range("D3").Value = 1.034603428466856
$ws.Range("E3").Value = 1.041843977526422
$ws.Range("F3").Value = 1.049991585345143
$ws.Range("I3").Value = 1.031845397144248
$ws.Range("J3").Value = 1.03854297930417
$ws.Range("K3").Value = 1.037213807026992
$ws.Range("L3").Value = 1.044435183564187
$ws.Range("M3").Value = 1.052561558337982
$ws.Range("N3").Value = 1.040017829457327
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034409532870776
$ws.Range("D4").Value = 1.035153865893396
$ws.Range("E4").Value = 1.042433560019326
$ws.Range("F4").Value = 1.050688914799804
$ws.Range("I4").Value = 1.031912716451113
$ws.Range("J4").Value = 1.038948022213424
$ws.Range("K4").Value = 1.037639972262671
$ws.Range("L4").Value = 1.04490136304331
$ws.Range("M4").Value = 1.05313629161197
$ws.Range("N4").Value = 1.040423447573971
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034677442923119
$ws.Range("D5").Value = 1.035385399810697
$ws.Range("E5").Value = 1.042681582109349
$ws.Range("F5").Value = 1.050982312760457
$ws.Range("I5").Value = 1.031940659507005
$ws.Range("J5").Value = 1.039118227234224
$ws.Range("K5").Value = 1.037819108780152
$ws.Range("L5").Value = 1.045097355761279
$ws.Range("M5").Value = 1.053378007139896
$ws.Range("N5").Value = 1.040593894305423
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034722430120167
$ws.Range("D6").Value = 1.035424282967569
$ws.Range("E6").Value = 1.042723235557573
$ws.Range("F6").Value = 1.051031589656362
$ws.Range("I6").Value = 1.0319453302633
$ws.Range("J6").Value = 1.039146800984524
$ws.Range("K6").Value = 1.037849185196282
$ws.Range("L6").Value = 1.045130264414714
$ws.Range("M6").Value = 1.053418597943431
$ws.Range("N6").Value = 1.040622508633726
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034413112434877
$ws.Range("D7").Value = 1.035156959151381
$ws.Range("E7").Value = 1.042436873468864
$ws.Range("F7").Value = 1.050692834253169
$ws.Range("I7").Value = 1.031913091234797
$ws.Range("J7").Value = 1.038950296796987
$ws.Range("K7").Value = 1.037642365986376
$ws.Range("L7").Value = 1.044903981865921
$ws.Range("M7").Value = 1.053139521043984
$ws.Range("N7").Value = 1.040425725387703
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033120786747304
$ws.Range("D8").Value = 1.03404069564328
$ws.Range("E8").Value = 1.041241300949814
$ws.Range("F8").Value = 1.049278936744648
$ws.Range("I8").Value = 1.031775300515261
$ws.Range("J8").Value = 1.038128306819214
$ws.Range("K8").Value = 1.036777699579077
$ws.Range("L8").Value = 1.043958252161453
$ws.Range("M8").Value = 1.051973856240851
$ws.Range("N8").Value = 1.039602568089879
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030845740066462
$ws.Range("D9").Value = 1.032077938098091
$ws.Range("E9").Value = 1.039139795265286
$ws.Range("F9").Value = 1.046795233225053
$ws.Range("I9").Value = 1.031521078875475
$ws.Range("J9").Value = 1.036677514414278
$ws.Range("K9").Value = 1.035253363452922
$ws.Range("L9").Value = 1.042292164531359
$ws.Range("M9").Value = 1.049922993085019
$ws.Range("N9").Value = 1.038149715393319
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029330554088519
$ws.Range("D10").Value = 1.030772347361021
$ws.Range("E10").Value = 1.037742376263886
$ws.Range("F10").Value = 1.045144700768771
$ws.Range("I10").Value = 1.031343949101374
$ws.Range("J10").Value = 1.035708790305116
$ws.Range("K10").Value = 1.034236728293135
$ws.Range("L10").Value = 1.041181761062283
$ws.Range("M10").Value = 1.048557941455348
$ws.Range("N10").Value = 1.037179615584808
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028674827173846
$ws.Range("D11").Value = 1.030207716212438
$ws.Range("E11").Value = 1.037138140423679
$ws.Range("F11").Value = 1.044431261809553
$ws.Range("I11").Value = 1.031265441400117
$ws.Range("J11").Value = 1.035288968441272
$ws.Range("K11").Value = 1.033796426391248
$ws.Range("L11").Value = 1.04070102982481
$ws.Range("M11").Value = 1.047967388133618
$ws.Range("N11").Value = 1.036759197525763
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028431315342692
$ws.Range("D12").Value = 1.029998092959099
$ws.Range("E12").Value = 1.036913829803085
$ws.Range("F12").Value = 1.044166447718897
$ws.Range("I12").Value = 1.031236008728759
$ws.Range("J12").Value = 1.035132974812548
$ws.Range("K12").Value = 1.033632865559378
$ws.Range("L12").Value = 1.040522477686558
$ws.Range("M12").Value = 1.047748109627683
$ws.Range("N12").Value = 1.036602982368193
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028483546990556
$ws.Range("D13").Value = 1.030043053050633
$ws.Range("E13").Value = 1.036961939319235
$ws.Range("F13").Value = 1.044223242680144
$ws.Range("I13").Value = 1.031242334412101
$ws.Range("J13").Value = 1.035166438367231
$ws.Range("K13").Value = 1.033667950493901
$ws.Range("L13").Value = 1.04056077713352
$ws.Range("M13").Value = 1.047795142011989
$ws.Range("N13").Value = 1.036636493444962
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028654697298678
$ws.Range("D14").Value = 1.030190386522963
$ws.Range("E14").Value = 1.037119596191376
$ws.Range("F14").Value = 1.044409368336362
$ws.Range("I14").Value = 1.031263014020701
$ws.Range("J14").Value = 1.035276075038469
$ws.Range("K14").Value = 1.033782906664804
$ws.Range("L14").Value = 1.040686270383161
$ws.Range("M14").Value = 1.047949260882257
$ws.Range("N14").Value = 1.03674628581285
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028760155866866
$ws.Range("D15").Value = 1.03028117759968
$ws.Range("E15").Value = 1.037216750977099
$ws.Range("F15").Value = 1.04452407154976
$ws.Range("I15").Value = 1.031275719457367
$ws.Range("J15").Value = 1.035343618794695
$ws.Range("K15").Value = 1.033853733242042
$ws.Range("L15").Value = 1.040763592639789
$ws.Range("M15").Value = 1.048044229133084
$ws.Range("N15").Value = 1.036813925488956
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029374080119301
$ws.Range("D16").Value = 1.030809834866776
$ws.Range("E16").Value = 1.037782495514173
$ws.Range("F16").Value = 1.045192075813757
$ws.Range("I16").Value = 1.031349121313849
$ws.Range("J16").Value = 1.035736645016319
$ws.Range("K16").Value = 1.034265947802061
$ws.Range("L16").Value = 1.04121366737402
$ws.Range("M16").Value = 1.048597145614895
$ws.Range("N16").Value = 1.037207509852896
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029759274920173
$ws.Range("D17").Value = 1.031141635245758
$ws.Range("E17").Value = 1.038137601873891
$ws.Range("F17").Value = 1.045611433082656
$ws.Range("I17").Value = 1.03139468024203
$ws.Range("J17").Value = 1.035983084632379
$ws.Range("K17").Value = 1.034524494955547
$ws.Range("L17").Value = 1.041496009495133
$ws.Range("M17").Value = 1.048944116050536
$ws.Range("N17").Value = 1.037454299441478
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029983987100134
$ws.Range("D18").Value = 1.031335236098392
$ws.Range("E18").Value = 1.038344811918443
$ws.Range("F18").Value = 1.045856157880084
$ws.Range("I18").Value = 1.031421079390954
$ws.Range("J18").Value = 1.036126794110152
$ws.Range("K18").Value = 1.034675292192868
$ws.Range("L18").Value = 1.041660702625134
$ws.Range("M18").Value = 1.049146548453513
$ws.Range("N18").Value = 1.037598213003193
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030060613995972
$ws.Range("D19").Value = 1.031401260372043
$ws.Range("E19").Value = 1.038415479157141
$ws.Range("F19").Value = 1.045939623178849
$ws.Range("I19").Value = 1.031430051190114
$ws.Range("J19").Value = 1.036175789442124
$ws.Range("K19").Value = 1.034726708614256
$ws.Range("L19").Value = 1.041716860018967
$ws.Range("M19").Value = 1.049215581201805
$ws.Range("N19").Value = 1.037647277914156
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029717943584365
$ws.Range("D20").Value = 1.031106029236329
$ws.Range("E20").Value = 1.038099493780954
$ws.Range("F20").Value = 1.045566427516454
$ws.Range("I20").Value = 1.031389810257189
$ws.Range("J20").Value = 1.035956647565762
$ws.Range("K20").Value = 1.034496756212427
$ws.Range("L20").Value = 1.041465716051049
$ws.Range("M20").Value = 1.048906884186307
$ws.Range("N20").Value = 1.037427824831193
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028604296294917
$ws.Range("D21").Value = 1.030146997556366
$ws.Range("E21").Value = 1.037073166590686
$ws.Range("F21").Value = 1.044354553749071
$ws.Range("I21").Value = 1.031256931880439
$ws.Range("J21").Value = 1.035243791231659
$ws.Range("K21").Value = 1.033749055293338
$ws.Range("L21").Value = 1.040649315388606
$ws.Range("M21").Value = 1.047903874521701
$ws.Range("N21").Value = 1.03671395615933
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027904416677535
$ws.Range("D22").Value = 1.029544629443886
$ws.Range("E22").Value = 1.036428623493886
$ws.Range("F22").Value = 1.043593693921712
$ws.Range("I22").Value = 1.031171815851079
$ws.Range("J22").Value = 1.034795283301462
$ws.Range("K22").Value = 1.033278870620364
$ws.Range("L22").Value = 1.04013608657559
$ws.Range("M22").Value = 1.047273701649892
$ws.Range("N22").Value = 1.036264811296428
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028275406125305
$ws.Range("D23").Value = 1.029863897725564
$ws.Range("E23").Value = 1.036770236636432
$ws.Range("F23").Value = 1.043996936233533
$ws.Range("I23").Value = 1.031217086158526
$ws.Range("J23").Value = 1.035033074701055
$ws.Range("K23").Value = 1.033528131312214
$ws.Range("L23").Value = 1.040408151563269
$ws.Range("M23").Value = 1.047607724455804
$ws.Range("N23").Value = 1.036502940387084
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029736619329068
$ws.Range("D24").Value = 1.031122117851446
$ws.Range("E24").Value = 1.038116712932899
$ws.Range("F24").Value = 1.045586763221334
$ws.Range("I24").Value = 1.031392011333234
$ws.Range("J24").Value = 1.035968593445063
$ws.Range("K24").Value = 1.034509290184106
$ws.Range("L24").Value = 1.041479404330546
$ws.Range("M24").Value = 1.048923707507348
$ws.Range("N24").Value = 1.037439787675014
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031433630284588
$ws.Range("D25").Value = 1.032584848890243
$ws.Range("E25").Value = 1.039682457246636
$ws.Range("F25").Value = 1.047436404895603
$ws.Range("I25").Value = 1.031588151395588
$ws.Range("J25").Value = 1.037052851835506
$ws.Range("K25").Value = 1.035647516880526
$ws.Range("L25").Value = 1.042722834852088
$ws.Range("M25").Value = 1.05045280777636
$ws.Range("N25").Value = 1.038525585836737
